$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update water-change-by-activities values (Blue/Green/Grey Water columns B:D, rows 2-11)
# Row 2 - CHEMICAL PROD
$ws.Range("B2").Value = -0.00001005974372425555
$ws.Range("C2").Value = -0.0002566238164192214
$ws.Range("D2").Value = -0.000006909280567235498

# Row 3 - COFFEE ESTATE
$ws.Range("B3").Value = -0.001494810133060298
$ws.Range("C3").Value = -0.02895324447661096
$ws.Range("D3").Value = -0.002092734311311517

# Row 4 - COOPERATIVES
$ws.Range("B4").Value = -0.0007449105419587987
$ws.Range("C4").Value = -140.8517357889305
$ws.Range("D4").Value = -0.006612574864448106

# Row 5 - ELECTRICITY PROD
$ws.Range("B5").Value = -0.0000002930676061829729
$ws.Range("C5").Value = -0.000003969572080819717
$ws.Range("D5").Value = -0.0000000972115061781409

# Row 6 - FERTILIZERS PROD
$ws.Range("B6").Value = -0.00001727742125229303
$ws.Range("C6").Value = -0.0004407465936751009
$ws.Range("D6").Value = -0.00001186655984328211

# Row 7 - INFORMAL
$ws.Range("B7").Value = -0.00001374847994739525
$ws.Range("C7").Value = -0.02771144274447579
$ws.Range("D7").Value = -0.0001220453300732061

# Row 8 - MANUFACTURING
$ws.Range("B8").Value = -0.0002988025216836832
$ws.Range("C8").Value = -0.01363722039241111
$ws.Range("D8").Value = -0.00007422279045954383

# Row 9 - PETROLEUM PROD
$ws.Range("B9").Value = 0.000002912342931704615
$ws.Range("C9").Value = 0.00008987740397969901
$ws.Range("D9").Value = 0.0001662855996187318

# Row 10 - PRIMARY
$ws.Range("B10").Value = -0.008530188227496183
$ws.Range("C10").Value = -0.1651569672576443
$ws.Range("D10").Value = -0.01193333795401941

# Row 11 - SERVICES
$ws.Range("B11").Value = -0.000001310376681384895
$ws.Range("D11").Value = -0.00002292652199287204
